$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Locate the "Método construtor " paragraph (robust text anchor via Find).
# ---------------------------------------------------------------------------
$findRng = $d.Content.Duplicate
$null = $findRng.Find.Execute("Método construtor ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

$origIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Start -eq $findRng.Start) {
        $origIndex = $i
        break
    }
}

# ---------------------------------------------------------------------------
# Insert 6 new blank paragraphs right before it - these will hold:
#   1) the "O Construtor..." paragraph   (moved up from later in the doc)
#   2) a fully empty trailing-style paragraph (moved up, matches <w:p/>)
#   3) new "Parâmetros são os valores..." paragraph
#   4) new blank paragraph
#   5) new "Parâmetros são características..." paragraph
#   6) new blank paragraph
# ---------------------------------------------------------------------------
$insertCount = 6
for ($i = 0; $i -lt $insertCount; $i++) {
    $idx = $origIndex + $i
    $d.Paragraphs($idx).Range.InsertParagraphBefore()
}

$p1 = $origIndex          # will hold "O Construtor..." text
$p2 = $origIndex + 1      # stays empty (becomes bare <w:p/>)
$p3 = $origIndex + 2      # "Parâmetros são os valores..."
$p4 = $origIndex + 3      # stays empty
$p5 = $origIndex + 4      # "Parâmetros são características..."
$p6 = $origIndex + 5      # stays empty
$heroIdx = $origIndex + 6 # the (still unmodified) "Método construtor " paragraph

# ---------------------------------------------------------------------------
# Fill in the moved / new paragraph text.
# ---------------------------------------------------------------------------
$d.Paragraphs($p1).Range.Text = "O Construtor é um método (a palavra construir é uma ação, ou seja, um método)"
$d.Paragraphs($p3).Range.Text = "Parâmetros são os valores essências para definir a característica"
$d.Paragraphs($p5).Range.Text = "Parâmetros são características obrigatórias, são predefinições ou filtros"

# ---------------------------------------------------------------------------
# Update the hero paragraph's own text, then split it into two runs with the
# pre-existing "_GoBack" bookmark sitting between them, exactly as in the
# target: "...classe he" | bookmark | "rde atributos e métodos de outra".
# ---------------------------------------------------------------------------
$heroPara = $d.Paragraphs($heroIdx)
$heroStart = $heroPara.Range.Start
$heroPara.Range.Text = "Herança permite que uma classe herde atributos e métodos de outra"

$splitPos = $heroStart + 33   # length of "Herança permite que uma classe he"
$bmRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# ---------------------------------------------------------------------------
# Remove the now-stale trailing content that followed the old hero paragraph:
#   - the blank paragraph right after it
#   - the old "O Construtor..." paragraph (now duplicated earlier in doc)
#   - the old trailing empty paragraph (also duplicated earlier)
# ---------------------------------------------------------------------------
$staleBlankIdx = $heroIdx + 1
$d.Paragraphs($staleBlankIdx).Range.Delete()

$staleConstrutorIdx = $heroIdx + 1
$d.Paragraphs($staleConstrutorIdx).Range.Delete()

$staleTrailIdx = $heroIdx + 1
$d.Paragraphs($staleTrailIdx).Range.Delete()

Write-Output "done"
